# Deploying to gh-pages from @ codeforIATI/codelists@91fa57d84c5da6f03203a179a8555ea2c1c6a4f1
#
# The sheet's columns E, F, G ("codeforiati:group-name", "codeforiati:category-name",
# "codeforiati:group-code" per the original header) were mislabelled: the
# group-code had ended up in the last column instead of right after the
# category-code. Fix: rotate the E/F/G values of every row (including the
# header) right by one -- the value that was in G moves to E, the value
# that was in E moves to F, and the value that was in F moves to G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$rng = $ws.Range("E1:G$lastRow")
$vals = $rng.Value2

$nrows = $vals.GetLength(0)
$newVals = New-Object 'object[,]' $nrows, 3

for ($i = 1; $i -le $nrows; $i++) {
    $eVal = $vals[$i, 1]
    $fVal = $vals[$i, 2]
    $gVal = $vals[$i, 3]

    $newVals[$i - 1, 0] = $gVal
    $newVals[$i - 1, 1] = $eVal
    $newVals[$i - 1, 2] = $fVal
}

# Force text format while writing so purely-numeric-looking strings (e.g.
# "110") stay stored as text, matching the original shared-string cell
# type, then restore the range's normal (default) style.
$rng.NumberFormat = "@"
$rng.Value = $newVals
$rng.Style = "Normal"

Write-Host "Rotated columns E/F/G for rows 1..$lastRow"
